# Daily attendance processing - 2025-12-12 16:57:57
# Reorders the "Recorded By" (column G) values: any non-"System" entries are
# reversed in order, and "System" entries (exact case) are moved to the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ([string]::IsNullOrEmpty($value)) {
        continue
    }

    $parts = $value -split ',\s*'

    $nonSystem = @()
    $systemCount = 0

    foreach ($part in $parts) {
        if ($part.Equals('System')) {
            $systemCount++
        } else {
            $nonSystem += $part
        }
    }

    if ($nonSystem.Count -gt 1) {
        $nonSystemRev = $nonSystem[($nonSystem.Count - 1)..0]
    } else {
        $nonSystemRev = $nonSystem
    }

    $newParts = @()
    $newParts += $nonSystemRev
    for ($i = 0; $i -lt $systemCount; $i++) { $newParts += 'System' }

    $newValue = $newParts -join ', '

    if ($newValue -ne $value) {
        $cell.Value2 = $newValue
    }
}
